$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the date in column C (rows 2-18) by 1 day (45171 -> 45172)
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 1
}
